$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Metadata" ---
$ws = $wb.Worksheets.Item(1)

# Remove the duplicated "Contact / No display for ContactDetail" row (old row 11).
$ws.Rows.Item(11).Delete()

# Version: 5.0.0 -> 6.0.0
$ws.Cells.Item(3, 2).Value = "6.0.0"

# Date: bump to new publish date
$ws.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty, now set to "Alvearie Team"
$ws.Cells.Item(9, 2).Value = "Alvearie Team"

# The remaining "Contact" row becomes "Jurisdiction / United States of America"
$ws.Cells.Item(10, 1).Value = "Jurisdiction"
$ws.Cells.Item(10, 2).Value = "United States of America"

# --- Sheet 2: "Include from Engagement Commu..." ---
$ws2 = $wb.Worksheets.Item(2)

# Fix typo in System URI: "reqeust" -> "request"
$ws2.Cells.Item(4, 2).Value = "http://ibm.com/fhir/cdm/CodeSystem/eng-communication-request-status-reason"
